# The deck's "last saved" date shown by the footer/date placeholder's
# datetimeFigureOut field needs to be refreshed from 4/13/2021 to
# 4/14/2021 on the slide master and every slide layout (this is the
# cached text PowerPoint stamps into the field when the date is not set
# to literally recompute itself at render time).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Type -ne 14) { continue }
        if ($shp.PlaceholderFormat.Type -ne 16) { continue }
        if (-not $shp.HasTextFrame) { continue }
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -ne $newText) {
            $tr.Text = $newText
        }
    }
}

$newDate = "4/14/2021"

# Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $newDate

# Every slide layout off the master.
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Update-DatePlaceholder $layout.Shapes $newDate
}
